$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 204, shifting existing data (rows 204+) down by two.
$ws.Rows.Item(204).Insert()
$ws.Rows.Item(204).Insert()

# Match the date-column number format used throughout column D.
$ws.Range("D204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D205").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 204 (Primera)
$ws.Range("A204").Value = 8
$ws.Range("B204").Value = "Terminal La Palmera de La Serena"
$ws.Range("C204").Value = "Coquimbo"
$ws.Range("D204").Value = 44524
$ws.Range("E204").Value = 4
$ws.Range("F204").Value = 100112009
$ws.Range("G204").Value = "Acelga"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 2600
$ws.Range("K204").Value = 550
$ws.Range("L204").Value = 600
$ws.Range("M204").Value = 575
$ws.Range("N204").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O204").Value = "Provincia del Elquí"
$ws.Range("P204").Value = 288
$ws.Range("Q204").Value = 2
$ws.Range("R204").Value = "Hortaliza"

# New row 205 (Segunda)
$ws.Range("A205").Value = 8
$ws.Range("B205").Value = "Terminal La Palmera de La Serena"
$ws.Range("C205").Value = "Coquimbo"
$ws.Range("D205").Value = 44524
$ws.Range("E205").Value = 4
$ws.Range("F205").Value = 100112009
$ws.Range("G205").Value = "Acelga"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Segunda"
$ws.Range("J205").Value = 1400
$ws.Range("K205").Value = 450
$ws.Range("L205").Value = 500
$ws.Range("M205").Value = 475
$ws.Range("N205").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O205").Value = "Provincia del Elquí"
$ws.Range("P205").Value = 238
$ws.Range("Q205").Value = 2
$ws.Range("R205").Value = "Hortaliza"
